$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add new summary block (rows 26-28), mirroring the header (row 2) and the
# --- Otsu+CC totals row (row 20), plus a new "CNN" results row ---

# Row 26: header row (same labels as row 2)
$ws.Range("C26").Value = "MSE"
$ws.Range("D26").Value = "SSIM"
$ws.Range("E26").Value = "NRMSE"
$ws.Range("F26").Value = "F-Score"
$ws.Range("G26").Value = "Custom F-Score"

# Row 27: CNN result row (new label "CNN", values equal to threshold-99.6 row 11)
$ws.Range("B27").Value = "CNN"
$ws.Range("C27").Value = 266.14846801757801
$ws.Range("D27").Value = 0.98157170825909001
$ws.Range("E27").Value = 0.89339673565194
$ws.Range("F27").Value = 0.59840038059424905
$ws.Range("G27").Value = 0.72044090374507597

# Row 28: Otsu+CC result row (duplicate of row 20)
$ws.Range("B28").Value = "Otsu+CC"
$ws.Range("C28").Value = 492.61874771118102
$ws.Range("D28").Value = 0.96178373089795299
$ws.Range("E28").Value = 0.80324402070840395
$ws.Range("F28").Value = 0.56775847140178604
$ws.Range("G28").Value = 0.61855112096831

# Copy the number formatting (Courier New font, left/center aligned) from the
# existing Otsu+CC row (row 20) onto the new numeric cells so they share the
# same style index instead of creating a brand-new style.
$ws.Range("C20:G20").Copy()
$ws.Range("C27:G27").PasteSpecial(-4122)
$ws.Range("C28:G28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match row heights used for the new block (15pt, custom height)
$ws.Rows.Item(26).RowHeight = 15
$ws.Rows.Item(27).RowHeight = 15
$ws.Rows.Item(28).RowHeight = 15

# --- Widen column K to fit the new, longer content ---
$ws.Columns.Item(11).ColumnWidth = 44.7109375

# --- Update the sheet view: scroll so row 7 is at the top, select K15 ---
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K15").Select()
